$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.877.37"
$ws.Range("E2").Value = "  -1.49%  "
$ws.Range("D3").Value = "2.043.06"
$ws.Range("E3").Value = "  -2.04%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.09"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.666"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.76"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.32"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.385"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0784"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.76%  "
$ws.Range("E12").Value = "  +1.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "16.21"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.68%  "
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.804"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.54"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +4.89%  "
$ws.Range("D17").Value = "2.045.65"
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("D18").Value = "36.876.30"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.74"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +13.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "75.19"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.36%  "
$ws.Range("D21").Value = "0.0₃0907"
$ws.Range("E21").Value = "  +6.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.39"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.66"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("E25").Value = "  -3.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.34"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +16.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.75"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.22"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.15"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.125"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("E31").Value = "  +3.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.72"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0618"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.44"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0880"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.70%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.21"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.73%  "
$ws.Range("E38").Value = "  -4.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.111"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +12.02%  "
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.75"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0223"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.13"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "96.57"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.83%  "
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.66"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +15.00%  "
$ws.Range("E47").Value = "  +4.01%  "
$ws.Range("D48").Value = "1.282.83"
$ws.Range("E48").Value = "  -3.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.89"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.73"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.90%  "
$ws.Range("D51").Value = "2.233.68"
$ws.Range("E51").Value = "  -1.93%  "
